{"js": "// Add a new paragraph \"A3\" at the end of the document body, after \"A2\".\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\"A3\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$lastPara = $d.Paragraphs.Last\n$endRange = $lastPara.Range\n$endRange.InsertParagraphAfter()\n\n$newRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)\n$newRange.Text = \"A3\"\n"}
